# Insert a new weekly price-report row at row 7 (Poroto granado, Vega Modelo
# de Temuco). This pushes the existing rows 7..46 down to 8..47, which is
# exactly the effect seen in the diff (every prior row's data shifts down by
# one row, and the dimension grows from R46 to R47).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly record.
$ws.Range("A7").Value = 10
$ws.Range("B7").Value = "Vega Modelo de Temuco"
$ws.Range("C7").Value = "La Araucanía"
$ws.Range("D7").Value = "2021-12-22"
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = 100112030
$ws.Range("G7").Value = "Poroto granado"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 110
$ws.Range("K7").Value = 45000
$ws.Range("L7").Value = 45000
$ws.Range("M7").Value = 45000
$ws.Range("N7").Value = "`$/saco 25 kilos"
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 1800
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = "Hortaliza"
